$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Bohemian Rhapsody" (row 2) measures were revised:
#   lines:   1112 -> 1592
#   seconds:  123 -> 166
# All other cells in the sheet (sums, D:I helper columns, the
# regression slope/intercept, and the correlation coefficient in B12)
# are formulas and will recompute automatically.
$ws.Range("B2").Value = 1592
$ws.Range("C2").Value = 166

# Force a full recalculation of every formula in the workbook.
$excel.CalculateFullRebuild()

# The saved view had the active cell on C3 when the file was last written.
$ws.Range("C3").Select()
